# feat: add 2022-Q1 data
#
# - Insert a new "2022-Q1" sheet (fund-level holdings) before the "总计" sheet.
# - Insert a new first data row into "总计" summarizing 2022-Q1, re-numbering
#   the index column for the rows that follow.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" worksheet, positioned right before "总计".
#
# Rather than inserting a brand-new blank sheet before "总计" (which would
# leave the existing "总计" sheet's old content/identity untouched), rename
# the current "总计" sheet to "2022-Q1", clear its old contents, and create a
# fresh "总计" sheet right after it. This keeps sheet order, relationship ids
# and the tab's identity aligned the same way the source edit produced them.
# ---------------------------------------------------------------------------
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"
$q1.Cells.Clear()

$total = $wb.Worksheets.Add($null, $q1)
$total.Name = "总计"

# Reuse an existing quarter sheet purely as a formatting template so the new
# sheet's styles (header border/bold/center, index-column style) match the
# rest of the workbook instead of inventing new ones.
$template = $wb.Worksheets.Item("2021-Q4")

# Header row (B1:H1) formatting.
$template.Range("B1:H1").Copy()
$q1.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats

# Index column (A2:A18) formatting.
$template.Range("A2").Copy()
$q1.Range("A2:A18").PasteSpecial(-4122)

# Data columns (B2:G18) formatting.
$template.Range("B2:G2").Copy()
$q1.Range("B2:G18").PasteSpecial(-4122)

# Rank column (H2:H18) formatting (plain/no style, like the template).
$template.Range("H2").Copy()
$q1.Range("H2:H18").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Columns B-G hold values that look numeric ("209.03", "002943", ...) but
# must stay text, matching the source data. Force text format before writing
# so Excel doesn't silently coerce them to numbers; the format is reverted
# to Normal afterwards (further down) once the literal text is locked in, so
# the cells end up back on the sheet's default (unstyled) look.
$dataBlock = $q1.Range("B2:G18")
$dataBlock.NumberFormat = "@"

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Data rows: index, code, name, scale, total position, position %, held
# value (亿元), position rank.
$rows = @(
    @(0,  "002943", "广发多因子灵活配置混合",             "209.03", "94.20", "3.96", "8.2776", 5),
    @(1,  "217024", "招商安盈债券",                         "35.05",  "20.20", "2.34", "0.8202", 2),
    @(2,  "003986", "申万菱信中证500指数优选增强A",         "23.29",  "92.19", "1.96", "0.4565", 3),
    @(3,  "014887", "招商安福1年定期开放债券",             "17.22",  "27.65", "2.09", "0.3599", 4),
    @(4,  "002616", "中银益利灵活配置混合A",               "5.90",   "29.85", "1.94", "0.1145", 9),
    @(5,  "006048", "长城中证500指数增强A",                 "4.45",   "92.64", "1.99", "0.0886", 9),
    @(6,  "160620", "鹏华中证A股资源产业指数（LOF）",       "2.77",   "94.14", "2.60", "0.0720", 1),
    @(7,  "007794", "申万菱信中证500指数优选增强C",         "3.53",   "92.19", "1.96", "0.0692", 3),
    @(8,  "007413", "长城中证500指数增强C",                 "1.72",   "92.64", "1.99", "0.0342", 9),
    @(9,  "005381", "泰康睿利量化多策略混合A",             "0.99",   "93.49", "2.33", "0.0231", 3),
    @(10, "011677", "中银睿丰回报混合型证券投资基金A",     "0.73",   "20.29", "1.96", "0.0143", 5),
    @(11, "005382", "泰康睿利量化多策略混合C",             "0.49",   "93.49", "2.33", "0.0114", 3),
    @(12, "002617", "中银益利灵活配置混合C",               "0.58",   "29.85", "1.94", "0.0113", 9),
    @(13, "006783", "红土创新中证500指数增强A",             "0.42",   "91.83", "2.14", "0.0090", 10),
    @(14, "006784", "红土创新中证500指数增强C",             "0.13",   "91.83", "2.14", "0.0028", 10),
    @(15, "515510", "嘉实中证500成长估值ETF",               "0.15",   "98.79", "1.22", "0.0018", 9),
    @(16, "011678", "中银睿丰回报混合型证券投资基金C",     "0.00",   "20.29", "1.96", $null,    5)
)

$r = 2
foreach ($row in $rows) {
    $q1.Range("A$r").Value = $row[0]
    $q1.Range("B$r").Value = $row[1]
    $q1.Range("C$r").Value = $row[2]
    $q1.Range("D$r").Value = $row[3]
    $q1.Range("E$r").Value = $row[4]
    $q1.Range("F$r").Value = $row[5]
    if ($row[6] -eq $null) {
        $g = $q1.Range("G$r")
        $g.NumberFormat = "General"
        $g.Value = 0
    } else {
        $q1.Range("G$r").Value = $row[6]
    }
    $q1.Range("H$r").Value = $row[7]
    $r++
}

# The "force text" number format has done its job (the literal strings are
# locked in as text); drop back to the sheet's normal/default style so these
# cells don't carry a lingering explicit format.
$dataBlock.Style = "Normal"

# ---------------------------------------------------------------------------
# 2) Populate the new "总计" sheet: the original summary table (quarters
#    2020-Q4 .. 2021-Q4) plus a new first row for 2022-Q1, index column
#    renumbered to match.
# ---------------------------------------------------------------------------
$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 17, 10.37),
    @("2021-Q4", 14, 2.31),
    @("2021-Q3", 13, 1.02),
    @("2021-Q2", 9,  2.29),
    @("2021-Q1", 1,  0.03),
    @("2020-Q4", 2,  1.82)
)

$r = 2
$idx = 0
foreach ($row in $totalRows) {
    $total.Range("A$r").Value = $idx
    $total.Range("B$r").Value = $row[0]
    $total.Range("C$r").Value = $row[1]
    $total.Range("D$r").Value = $row[2]
    $r++
    $idx++
}

# Formatting: header (B1:D1) bold/centered/bordered, index column (A2:A7)
# matching style, reusing the "2022-Q1" sheet (itself cloned from a quarter
# sheet) as the style template so no new style entries are invented.
$q1.Range("B1:D1").Copy()
$total.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2").Copy()
$total.Range("A2:A7").PasteSpecial(-4122)
$excel.CutCopyMode = $false
